$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 298.46155
$ws.Range("I19").Value = 173.75
$ws.Range("J19").Value = 353.8889
$ws.Range("K19").Value = 173.75
$ws.Range("L19").Value = 353.8889
$ws.Range("M19").Value = 1.25
$ws.Range("N19").Value = -703.8888999999999
$ws.Range("H111").Value = 900
$ws.Range("I111").Value = 700
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 2100
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 967
$ws.Range("N111").Value = -9134
$ws.Range("H127").Value = 2141.5
$ws.Range("I127").Value = 1049
$ws.Range("J127").Value = 2439.4546
$ws.Range("K127").Value = 3147
$ws.Range("L127").Value = 7318.3638
$ws.Range("M127").Value = 1813
$ws.Range("N127").Value = -17238.3638
$ws.Range("H129").Value = 1329.9286
$ws.Range("I129").Value = 4098.5
$ws.Range("J129").Value = 868.5
$ws.Range("K129").Value = 12295.5
$ws.Range("L129").Value = 2605.5
$ws.Range("M129").Value = -7295.5
$ws.Range("H132").Value = 21991.553
$ws.Range("I132").Value = 3062.842
$ws.Range("J132").Value = 101912.78
$ws.Range("K132").Value = 9188.526
$ws.Range("L132").Value = 305738.34
$ws.Range("M132").Value = -6658.526
$ws.Range("N132").Value = -310798.34
$ws.Range("H137").Value = 11078.131
$ws.Range("I137").Value = 11352.667
$ws.Range("J137").Value = 10778.637
$ws.Range("K137").Value = 34058.001
$ws.Range("L137").Value = 32335.911
$ws.Range("M137").Value = -31508.001
$ws.Range("N137").Value = -37435.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9936.770500000001
$ws.Range("I32").Value = 9672.394
$ws.Range("J32").Value = 11177.308
$ws.Range("K32").Value = 9672.394
$ws.Range("L32").Value = 11177.308
$ws.Range("M32").Value = -9385.394
$ws.Range("N32").Value = -11751.308
$ws.Range("H61").Value = 2093.6428
$ws.Range("I61").Value = 1605.8
$ws.Range("J61").Value = 2811.0588
$ws.Range("K61").Value = 1605.8
$ws.Range("L61").Value = 2811.0588
$ws.Range("M61").Value = -1393.8
$ws.Range("N61").Value = -3235.0588
$ws.Range("H74").Value = 1646.079
$ws.Range("I74").Value = 1337.5555
$ws.Range("J74").Value = 7199.5
$ws.Range("K74").Value = 1337.5555
$ws.Range("L74").Value = 7199.5
$ws.Range("M74").Value = -463.5554999999999
$ws.Range("N74").Value = -8947.5
$ws.Range("H77").Value = 1646.079
$ws.Range("I77").Value = 1337.5555
$ws.Range("J77").Value = 7199.5
$ws.Range("K77").Value = 6687.7775
$ws.Range("L77").Value = 35997.5
$ws.Range("M77").Value = -2319.7775
$ws.Range("N77").Value = -44733.5
$ws.Range("H110").Value = 4308.6
$ws.Range("I110").Value = 2012.2858
$ws.Range("J110").Value = 9666.666999999999
$ws.Range("K110").Value = 2012.2858
$ws.Range("L110").Value = 9666.666999999999
$ws.Range("M110").Value = 32.71419999999989
$ws.Range("H122").Value = 2134.6667
$ws.Range("I122").Value = 2064
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 6192
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -3742
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 13160521
$ws.Range("I132").Value = 19232346
$ws.Range("J132").Value = 4901.75
$ws.Range("K132").Value = 57697038
$ws.Range("L132").Value = 14705.25
$ws.Range("M132").Value = -57694508
$ws.Range("N132").Value = -19765.25
$ws.Range("H136").Value = 2093.6428
$ws.Range("I136").Value = 1605.8
$ws.Range("J136").Value = 2811.0588
$ws.Range("K136").Value = 4817.4
$ws.Range("L136").Value = 8433.1764
$ws.Range("M136").Value = -2267.4
$ws.Range("N136").Value = -13533.1764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5072.143
$ws.Range("I86").Value = 4301.2
$ws.Range("J86").Value = 6999.5
$ws.Range("K86").Value = 4301.2
$ws.Range("L86").Value = 6999.5
$ws.Range("M86").Value = -3178.2
$ws.Range("N86").Value = -9245.5
$ws.Range("H89").Value = 5072.143
$ws.Range("I89").Value = 4301.2
$ws.Range("J89").Value = 6999.5
$ws.Range("K89").Value = 21506
$ws.Range("L89").Value = 34997.5
$ws.Range("M89").Value = -15890
$ws.Range("N89").Value = -46229.5
$ws.Range("H132").Value = 52000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 52000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 52000
$ws.Range("N132").Value = -62120
$ws.Range("H134").Value = 1750.6296
$ws.Range("I134").Value = 1225.909
$ws.Range("J134").Value = 4059.4
$ws.Range("K134").Value = 3677.727
$ws.Range("L134").Value = 12178.2
$ws.Range("M134").Value = -1142.727
$ws.Range("N134").Value = -17248.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10205346
$ws.Range("I58").Value = 753.6667
$ws.Range("J58").Value = 31252320
$ws.Range("K58").Value = 753.6667
$ws.Range("L58").Value = 31252320
$ws.Range("M58").Value = -550.6667
$ws.Range("N58").Value = -31252726
$ws.Range("H94").Value = 975.75
$ws.Range("I94").Value = 356.5
$ws.Range("J94").Value = 1099.6
$ws.Range("K94").Value = 356.5
$ws.Range("L94").Value = 1099.6
$ws.Range("M94").Value = 94.5
$ws.Range("N94").Value = -2001.6
$ws.Range("H99").Value = 3126.5
$ws.Range("I99").Value = 3126.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3126.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1628.5
$ws.Range("H122").Value = 120841.1
$ws.Range("I122").Value = 134179
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 402537
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -400087
$ws.Range("N122").Value = -7300
$ws.Range("H126").Value = 3126.5
$ws.Range("I126").Value = 3126.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9379.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6909.5
$ws.Range("H136").Value = 10205346
$ws.Range("I136").Value = 753.6667
$ws.Range("J136").Value = 31252320
$ws.Range("K136").Value = 2261.0001
$ws.Range("L136").Value = 93756960
$ws.Range("M136").Value = 288.9998999999998
$ws.Range("N136").Value = -93762060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 209.28572
$ws.Range("I12").Value = 282.375
$ws.Range("J12").Value = 164.3077
$ws.Range("K12").Value = 847.125
$ws.Range("L12").Value = 492.9231
$ws.Range("M12").Value = -674.125
$ws.Range("N12").Value = -838.9231
$ws.Range("H23").Value = 581.2727
$ws.Range("I23").Value = 555.2857
$ws.Range("J23").Value = 593.4
$ws.Range("K23").Value = 1665.8571
$ws.Range("L23").Value = 1780.2
$ws.Range("M23").Value = -1430.8571
$ws.Range("N23").Value = -2250.2
$ws.Range("H68").Value = 1232.9362
$ws.Range("I68").Value = 1116.9
$ws.Range("J68").Value = 1264.2972
$ws.Range("K68").Value = 3350.7
$ws.Range("L68").Value = 3792.8916
$ws.Range("M68").Value = -2539.7
$ws.Range("N68").Value = -5414.8916
$ws.Range("H71").Value = 1232.9362
$ws.Range("I71").Value = 1116.9
$ws.Range("J71").Value = 1264.2972
$ws.Range("K71").Value = 10052.1
$ws.Range("L71").Value = 11378.6748
$ws.Range("M71").Value = -5996.1
$ws.Range("N71").Value = -19490.6748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4212.8335
$ws.Range("I122").Value = 3853.8
$ws.Range("J122").Value = 6008
$ws.Range("K122").Value = 11561.4
$ws.Range("L122").Value = 18024
$ws.Range("M122").Value = -9111.400000000001
$ws.Range("N122").Value = -22924
$ws.Range("H132").Value = 2446.6155
$ws.Range("I132").Value = 1886.7142
$ws.Range("J132").Value = 4798.2
$ws.Range("K132").Value = 5660.142599999999
$ws.Range("L132").Value = 14394.6
$ws.Range("M132").Value = -3130.142599999999
$ws.Range("N132").Value = -19454.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 52634264
$ws.Range("I7").Value = 62502130
$ws.Range("J7").Value = 5665
$ws.Range("K7").Value = 62502130
$ws.Range("L7").Value = 5665
$ws.Range("M7").Value = -62502018
$ws.Range("N7").Value = -5889
$ws.Range("H40").Value = 5542.6313
$ws.Range("I40").Value = 5346.769
$ws.Range("J40").Value = 5967
$ws.Range("K40").Value = 5346.769
$ws.Range("L40").Value = 5967
$ws.Range("M40").Value = -5210.769
$ws.Range("N40").Value = -6239
$ws.Range("H46").Value = 4608.8696
$ws.Range("I46").Value = 871.5714
$ws.Range("J46").Value = 10422.444
$ws.Range("K46").Value = 871.5714
$ws.Range("L46").Value = 10422.444
$ws.Range("M46").Value = -683.5714
$ws.Range("N46").Value = -10798.444
$ws.Range("H122").Value = 2337.4375
$ws.Range("I122").Value = 2335.6428
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 7006.928400000001
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -4556.928400000001
$ws.Range("N122").Value = -11950
$ws.Range("H126").Value = 52634264
$ws.Range("I126").Value = 62502130
$ws.Range("J126").Value = 5665
$ws.Range("K126").Value = 187506390
$ws.Range("L126").Value = 16995
$ws.Range("M126").Value = -187503920
$ws.Range("N126").Value = -21935
$ws.Range("H132").Value = 2846.054
$ws.Range("I132").Value = 2178
$ws.Range("J132").Value = 4649.8
$ws.Range("K132").Value = 6534
$ws.Range("L132").Value = 13949.4
$ws.Range("M132").Value = -4004
$ws.Range("N132").Value = -19009.4
$ws.Range("H136").Value = 1415.8043
$ws.Range("I136").Value = 1009.13513
$ws.Range("J136").Value = 3087.6667
$ws.Range("K136").Value = 3027.40539
$ws.Range("L136").Value = 9263.000100000001
$ws.Range("M136").Value = -477.4053899999999
$ws.Range("N136").Value = -14363.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H107").Value = 6250624.5
$ws.Range("I107").Value = 565.75
$ws.Range("J107").Value = 25000800
$ws.Range("K107").Value = 1697.25
$ws.Range("L107").Value = 75002400
$ws.Range("M107").Value = 222.75
$ws.Range("N107").Value = -75006240
$ws.Range("H122").Value = 1091
$ws.Range("I122").Value = 1122.4445
$ws.Range("J122").Value = 996.6667
$ws.Range("K122").Value = 3367.3335
$ws.Range("L122").Value = 2990.0001
$ws.Range("M122").Value = -917.3335000000002
$ws.Range("N122").Value = -7890.0001
$ws.Range("H126").Value = 4903711.5
$ws.Range("I126").Value = 4903711.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14711134.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14708664.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1518.4857
$ws.Range("I132").Value = 1024.9615
$ws.Range("J132").Value = 2944.2222
$ws.Range("K132").Value = 3074.8845
$ws.Range("L132").Value = 8832.6666
$ws.Range("M132").Value = -544.8844999999997
$ws.Range("N132").Value = -13892.6666
$ws.Range("H136").Value = 200913.77
$ws.Range("I136").Value = 244611.83
$ws.Range("J136").Value = 1844.7778
$ws.Range("K136").Value = 733835.49
$ws.Range("L136").Value = 5534.3334
$ws.Range("M136").Value = -731285.49
$ws.Range("N136").Value = -10634.3334
